$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.756.84'
$ws.Range('E2').Value = '  -0.07%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.530.94'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.55'
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.15'
$ws.Range('E6').Value = '  +0.88%  '

# Row 7
$ws.Range('E7').Value = '  -0.39%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('E9').Value = '  -5.56%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.646'
$ws.Range('E10').Value = '  -2.20%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.24'
$ws.Range('E11').Value = '  -0.30%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000302'
$ws.Range('E12').Value = '  -1.46%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.46'
$ws.Range('E13').Value = '  -1.36%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.089.41'
$ws.Range('E14').Value = '  +0.79%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '592.98'
$ws.Range('E15').Value = '  -3.93%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.873.41'
$ws.Range('E16').Value = '  -0.05%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.73'
$ws.Range('E17').Value = '  +0.50%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.99'
$ws.Range('E18').Value = '  +0.87%  '

# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.537.60'
$ws.Range('E19').Value = '  +0.64%  '

# Row 20
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.122'
$ws.Range('E20').Value = '  +1.79%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.982'
$ws.Range('E21').Value = '  -0.69%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.77'
$ws.Range('E22').Value = '  -0.22%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '102.57'
$ws.Range('E23').Value = '  -2.82%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.12'
$ws.Range('E24').Value = '  +1.27%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.64'
$ws.Range('E25').Value = '  +0.14%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.04'
$ws.Range('E26').Value = '  +0.26%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.76'
$ws.Range('E27').Value = '  -2.10%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.52'
$ws.Range('E28').Value = '  -3.84%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.13'
$ws.Range('E29').Value = '  -3.23%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.03'
$ws.Range('E30').Value = '  -1.69%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.21'
$ws.Range('E31').Value = '  -1.05%  '

# Row 32
$ws.Range('E32').Value = '  -2.37%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('E33').Value = '  +0.00%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.37'
$ws.Range('E34').Value = '  -1.31%  '

# Row 35
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.20'
$ws.Range('E35').Value = '  +3.86%  '

# Row 36
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.810.95'
$ws.Range('E36').Value = '  +2.15%  '

# Row 37
$ws.Range('E37').Value = '  +0.07%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0809'
$ws.Range('E38').Value = '  +1.73%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '510.16'
$ws.Range('E39').Value = '  -2.45%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.391'
$ws.Range('E40').Value = '  +0.24%  '

# Row 41
$ws.Range('E41').Value = '  -0.20%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.57'
$ws.Range('E42').Value = '  -0.19%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.134'
$ws.Range('E43').Value = '  -2.69%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0448'
$ws.Range('E44').Value = '  -3.21%  '

# Row 45
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.35'
$ws.Range('E45').Value = '  +0.94%  '

# Row 46
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.139'
$ws.Range('E46').Value = '  -0.90%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  -1.77%  '

# Row 48
$ws.Range('E48').Value = '  +0.08%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.46'
$ws.Range('E49').Value = '  -3.21%  '

# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.32'
$ws.Range('E50').Value = '  +1.87%  '

# Row 51
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000244'
$ws.Range('E51').Value = '  +2.90%  '
